# Add a link to the CDC report.
#
# A new row is inserted just above the old row 17 ("CALCULATIONS" header),
# pushing that header and the calculation rows below it down by one row.
# The newly created row gets the hyperlink text and an actual hyperlink
# pointing at the CDC COVID data tracker report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, empty row before what used to be row 17, shifting
# rows 17-20 down to 18-21.
$ws.Rows("17:17").Insert()

# Put the link caption text in the new row.
$ws.Range("A16").Value = "  - see CDC report "

# Turn that cell into a hyperlink pointing at the CDC report. This also
# applies the built-in "Hyperlink" cell style (underlined, themed font).
$ws.Hyperlinks.Add(
    $ws.Range("A16"),
    "https://www.cdc.gov/coronavirus/2019-ncov/covid-data/covidview/index.html"
)

# Reflect the author's final cursor position on the new row.
[void]$ws.Range("A16").Select()
